$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row touched: H86, I86, J86, K86, L86, M86, N86
$ws.Range("H86").Value = 45457464
$ws.Range("I86").Value = 2312
$ws.Range("J86").Value = 71431830
$ws.Range("K86").Value = 2312
$ws.Range("L86").Value = 71431830
$ws.Range("M86").Value = -1189
$ws.Range("N86").Value = -71434076

# row touched: H89, I89, J89, K89, L89, M89, N89
$ws.Range("H89").Value = 45457464
$ws.Range("I89").Value = 2312
$ws.Range("J89").Value = 71431830
$ws.Range("K89").Value = 11560
$ws.Range("L89").Value = 357159150
$ws.Range("M89").Value = -5944
$ws.Range("N89").Value = -357170382

# row touched: H103, I103, J103, K103, L103, M103, N103
$ws.Range("H103").Value = 743.6667
$ws.Range("I103").Value = 826.7692
$ws.Range("J103").Value = 645.4545000000001
$ws.Range("K103").Value = 2480.3076
$ws.Range("L103").Value = 1936.3635
$ws.Range("M103").Value = -1894.3076
$ws.Range("N103").Value = -3108.3635

# row touched: H129, I129, J129, K129, L129, M129, N129
$ws.Range("H129").Value = 51852640
$ws.Range("I129").Value = 125000420
$ws.Range("J129").Value = 3087451.2
$ws.Range("K129").Value = 375001260
$ws.Range("L129").Value = 9262353.600000001
$ws.Range("M129").Value = -374996260
$ws.Range("N129").Value = -9272353.600000001

# row touched: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 2703.2917
$ws.Range("I132").Value = 2818.3635
$ws.Range("J132").Value = 1437.5
$ws.Range("K132").Value = 8455.0905
$ws.Range("L132").Value = 4312.5
$ws.Range("M132").Value = -5925.0905
$ws.Range("N132").Value = -9372.5

# row touched: H137, I137, J137, K137, L137, M137, N137
$ws.Range("H137").Value = 11112499
$ws.Range("I137").Value = 1411.4166
$ws.Range("J137").Value = 55556850
$ws.Range("K137").Value = 4234.2498
$ws.Range("L137").Value = 166670550
$ws.Range("M137").Value = -1684.2498
$ws.Range("N137").Value = -166675650

$ws = $wb.Worksheets.Item("ARM")
# row touched: H32, I32, J32, K32, L32, M32, N32
$ws.Range("H32").Value = 18588.984
$ws.Range("I32").Value = 18943.066
$ws.Range("J32").Value = 14340
$ws.Range("K32").Value = 18943.066
$ws.Range("L32").Value = 14340
$ws.Range("M32").Value = -18656.066
$ws.Range("N32").Value = -14914

# row touched: H48, I48, J48, K48, L48, N48
$ws.Range("H48").Value = 120000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 120000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 120000
$ws.Range("N48").Value = -120768

# row touched: H88, I88, J88, K88, L88, M88, N88
$ws.Range("H88").Value = 7862
$ws.Range("I88").Value = 3998
$ws.Range("J88").Value = 9150
$ws.Range("K88").Value = 3998
$ws.Range("L88").Value = 9150
$ws.Range("M88").Value = -3592
$ws.Range("N88").Value = -9962

# row touched: H91, I91, J91, K91, L91, M91, N91
$ws.Range("H91").Value = 7862
$ws.Range("I91").Value = 3998
$ws.Range("J91").Value = 9150
$ws.Range("K91").Value = 3998
$ws.Range("L91").Value = 9150
$ws.Range("M91").Value = -2594
$ws.Range("N91").Value = -11958

$ws = $wb.Worksheets.Item("BSM")
# row touched: H52, I52, J52, K52, L52, N52
$ws.Range("H52").Value = 14800
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 14800
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 14800
$ws.Range("N52").Value = -15326

# row touched: H86, I86, J86, K86, L86, M86, N86
$ws.Range("H86").Value = 3122.5
$ws.Range("I86").Value = 2690.1
$ws.Range("J86").Value = 4203.5
$ws.Range("K86").Value = 2690.1
$ws.Range("L86").Value = 4203.5
$ws.Range("M86").Value = -1567.1
$ws.Range("N86").Value = -6449.5

# row touched: H89, I89, J89, K89, L89, M89, N89
$ws.Range("H89").Value = 3122.5
$ws.Range("I89").Value = 2690.1
$ws.Range("J89").Value = 4203.5
$ws.Range("K89").Value = 13450.5
$ws.Range("L89").Value = 21017.5
$ws.Range("M89").Value = -7834.5
$ws.Range("N89").Value = -32249.5

# row touched: H99, I99, J99, K99, L99, M99, N99
$ws.Range("H99").Value = 55557610
$ws.Range("I99").Value = 66668692
$ws.Range("J99").Value = 2200
$ws.Range("K99").Value = 66668692
$ws.Range("L99").Value = 2200
$ws.Range("M99").Value = -66667194
$ws.Range("N99").Value = -5196

# row touched: H121, I121, J121, K121, L121, N121
$ws.Range("H121").Value = 14800
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 14800
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 14800
$ws.Range("N121").Value = -18294

$ws = $wb.Worksheets.Item("CRP")
# row touched: H3, I3, J3, K3, L3, M3, N3
$ws.Range("H3").Value = 900
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = -887
$ws.Range("N3").Value = -1026

# row touched: H4, I4, J4, K4, L4, M4, N4
$ws.Range("H4").Value = 11000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 11000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 11000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -11224

# row touched: H31, I31, J31, K31, L31, M31, N31
$ws.Range("H31").Value = 2007.8948
$ws.Range("I31").Value = 1544.8182
$ws.Range("J31").Value = 2644.625
$ws.Range("K31").Value = 1544.8182
$ws.Range("L31").Value = 2644.625
$ws.Range("M31").Value = -1249.8182
$ws.Range("N31").Value = -3234.625

# row touched: H34, I34, J34, K34, L34, M34, N34
$ws.Range("H34").Value = 2007.8948
$ws.Range("I34").Value = 1544.8182
$ws.Range("J34").Value = 2644.625
$ws.Range("K34").Value = 1544.8182
$ws.Range("L34").Value = 2644.625
$ws.Range("M34").Value = -1342.8182
$ws.Range("N34").Value = -3048.625

# row touched: H48, I48, J48, K48, L48, M48, N48
$ws.Range("H48").Value = 5348.6665
$ws.Range("I48").Value = 2046
$ws.Range("J48").Value = 7000
$ws.Range("K48").Value = 2046
$ws.Range("L48").Value = 7000
$ws.Range("M48").Value = -1570
$ws.Range("N48").Value = -7952

# row touched: H109, I109, J109, K109, L109, N109
$ws.Range("H109").Value = 23285
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 23285
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 23285
$ws.Range("N109").Value = -25365

# row touched: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 3493.2354
$ws.Range("I132").Value = 2793.6
$ws.Range("J132").Value = 4492.7144
$ws.Range("K132").Value = 8380.799999999999
$ws.Range("L132").Value = 13478.1432
$ws.Range("M132").Value = -5850.799999999999
$ws.Range("N132").Value = -18538.1432

# row touched: H133, I133, J133, K133, L133, M133, N133
$ws.Range("H133").Value = 70775.336
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 70775.336
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 70775.336
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -75835.336

$ws = $wb.Worksheets.Item("CUL")
# row touched: H2, I2, J2, K2, L2, M2, N2
$ws.Range("H2").Value = 204.25
$ws.Range("I2").Value = 158.11111
$ws.Range("J2").Value = 342.66666
$ws.Range("K2").Value = 948.66666
$ws.Range("L2").Value = 2055.99996
$ws.Range("M2").Value = -835.66666
$ws.Range("N2").Value = -2281.99996

# row touched: H4, I4, J4, K4, L4, M4, N4
$ws.Range("H4").Value = 2066.2222
$ws.Range("I4").Value = 1074.5
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 3223.5
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = -3111.5
$ws.Range("N4").Value = -30224

# row touched: H107, I107, J107, K107, L107, M107, N107
$ws.Range("H107").Value = 914
$ws.Range("I107").Value = 267.2857
$ws.Range("J107").Value = 1262.2307
$ws.Range("K107").Value = 801.8571000000001
$ws.Range("L107").Value = 3786.6921
$ws.Range("M107").Value = 1118.1429
$ws.Range("N107").Value = -7626.6921

# row touched: H122, I122, J122, K122, L122, M122
$ws.Range("H122").Value = 371.8889
$ws.Range("I122").Value = 371.8889
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3347.0001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -897.0000999999997

# row touched: H131, I131, J131, K131, L131, M131, N131
$ws.Range("H131").Value = 16149.928
$ws.Range("I131").Value = 72620.36
$ws.Range("J131").Value = 1775.6364
$ws.Range("K131").Value = 217861.08
$ws.Range("L131").Value = 5326.9092
$ws.Range("M131").Value = -212821.08
$ws.Range("N131").Value = -15406.9092

$ws = $wb.Worksheets.Item("LTW")
# row touched: H16, I16, J16, K16, L16, M16, N16
$ws.Range("H16").Value = 16428.572
$ws.Range("I16").Value = 19066.666
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 19066.666
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = -18896.666
$ws.Range("N16").Value = -940

# row touched: H100, I100, J100, K100, L100, M100, N100
$ws.Range("H100").Value = 7938226.5
$ws.Range("I100").Value = 18519836
$ws.Range("J100").Value = 2020
$ws.Range("K100").Value = 18519836
$ws.Range("L100").Value = 2020
$ws.Range("M100").Value = -18519295
$ws.Range("N100").Value = -3102

$ws = $wb.Worksheets.Item("WVR")
# row touched: H109, I109, J109, K109, L109, N109
$ws.Range("H109").Value = 24571.428
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 24571.428
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 24571.428
$ws.Range("N109").Value = -27345.428

# row touched: H122, I122, J122, K122, L122, M122, N122
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 4500
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -9400

# row touched: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 2597.0715
$ws.Range("I132").Value = 1450
$ws.Range("J132").Value = 3055.9
$ws.Range("K132").Value = 4350
$ws.Range("L132").Value = 9167.700000000001
$ws.Range("M132").Value = -1820
$ws.Range("N132").Value = -14227.7

# row touched: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 2263.0454
$ws.Range("I136").Value = 2338.423
$ws.Range("J136").Value = 2154.1667
$ws.Range("K136").Value = 7015.268999999999
$ws.Range("L136").Value = 6462.500100000001
$ws.Range("M136").Value = -4465.268999999999
$ws.Range("N136").Value = -11562.5001
